$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 149, pushing the existing rows 149-249
# down to 151-251 (so the old 248/249 end up as the new 250/251).
$ws.Rows.Item(149).Insert()
$ws.Rows.Item(149).Insert()

# New row 149 - Primera quality entry for the new date
$ws.Cells.Item(149,1).Value = 1
$ws.Cells.Item(149,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(149,3).Value = "Arica y Parinacota"
$ws.Cells.Item(149,4).Value = 44603
$ws.Cells.Item(149,5).Value = 15
$ws.Cells.Item(149,6).Value = 100112043
$ws.Cells.Item(149,7).Value = "Pepino ensalada"
$ws.Cells.Item(149,8).Value = "Sin especificar"
$ws.Cells.Item(149,9).Value = "Primera"
$ws.Cells.Item(149,10).Value = 130
$ws.Cells.Item(149,11).Value = 10000
$ws.Cells.Item(149,12).Value = 11000
$ws.Cells.Item(149,13).Value = 10500
$ws.Cells.Item(149,14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(149,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(149,16).Value = 150
$ws.Cells.Item(149,17).Value = 70
$ws.Cells.Item(149,18).Value = "Hortaliza"

# New row 150 - Segunda quality entry for the new date
$ws.Cells.Item(150,1).Value = 1
$ws.Cells.Item(150,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(150,3).Value = "Arica y Parinacota"
$ws.Cells.Item(150,4).Value = 44603
$ws.Cells.Item(150,5).Value = 15
$ws.Cells.Item(150,6).Value = 100112043
$ws.Cells.Item(150,7).Value = "Pepino ensalada"
$ws.Cells.Item(150,8).Value = "Sin especificar"
$ws.Cells.Item(150,9).Value = "Segunda"
$ws.Cells.Item(150,10).Value = 160
$ws.Cells.Item(150,11).Value = 8000
$ws.Cells.Item(150,12).Value = 9000
$ws.Cells.Item(150,13).Value = 8500
$ws.Cells.Item(150,14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(150,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(150,16).Value = 85
$ws.Cells.Item(150,17).Value = 100
$ws.Cells.Item(150,18).Value = "Hortaliza"
